$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 1.17
$ws.Range("L2").Value = 1.83
$ws.Range("M2").Value = 1.83
$ws.Range("AA4").Value = 21
$ws.Range("AD4").Value = 12.5
$ws.Range("AE4").Value = 7.7
$ws.Range("AF4").Value = 12.5
$ws.Range("H4").Value = 7.9
$ws.Range("N4").Value = 1.24
$ws.Range("O4").Value = 3.3
$ws.Range("Z4").Value = 24
$ws.Range("AA5").Value = 9.75
$ws.Range("AB5").Value = 18.5
$ws.Range("AC5").Value = 75
$ws.Range("AD5").Value = 8.25
$ws.Range("AE5").Value = 7.3
$ws.Range("AF5").Value = 8.5
$ws.Range("AG5").Value = 9
$ws.Range("AI5").Value = 23
$ws.Range("AJ5").Value = 500
$ws.Range("G5").Value = 7.1
$ws.Range("H5").Value = 4.75
$ws.Range("I5").Value = 1.36
$ws.Range("N5").Value = 1.5
$ws.Range("O5").Value = 2.25
$ws.Range("R5").Value = 1.75
$ws.Range("S5").Value = 1.85
$ws.Range("U5").Value = 50
$ws.Range("V5").Value = 22
$ws.Range("W5").Value = 175
$ws.Range("X5").Value = 75
$ws.Range("Y5").Value = 65
$ws.Range("Z5").Value = 15.5
$ws.Range("AA6").Value = 14
$ws.Range("AB6").Value = 24
$ws.Range("AE6").Value = 150
$ws.Range("AF6").Value = 45
$ws.Range("AG6").Value = 600
$ws.Range("AJ6").Value = 500
$ws.Range("H6").Value = 6.2
$ws.Range("I6").Value = 13.5
$ws.Range("N6").Value = 1.32
$ws.Range("O6").Value = 2.82
$ws.Range("R6").Value = 1.82
$ws.Range("S6").Value = 1.78
$ws.Range("T6").Value = 10.5
$ws.Range("U6").Value = 7.4
$ws.Range("V6").Value = 9.75
$ws.Range("Y6").Value = 25
$ws.Range("Z6").Value = 21
$ws.Range("AC7").Value = 350
$ws.Range("AD7").Value = 32
$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 16.5
$ws.Range("R7").Value = 2.9
$ws.Range("AD8").Value = 8.25
$ws.Range("AG8").Value = 18
$ws.Range("J8").Value = 1.04
$ws.Range("K8").Value = 9
$ws.Range("Z8").Value = 9
$ws.Range("J14").Value = 1.1
$ws.Range("K14").Value = 7
$ws.Range("AA15").Value = 6
$ws.Range("AB15").Value = 15
$ws.Range("AE15").Value = 9
$ws.Range("AH15").Value = 17
$ws.Range("AI15").Value = 29
$ws.Range("AJ15").Value = 351
$ws.Range("K15").Value = 8.5
$ws.Range("T15").Value = 10
$ws.Range("Z15").Value = 8.5
$ws.Range("N16").Value = 1.8
$ws.Range("O16").Value = 2
$ws.Range("N18").Value = 2.05
$ws.Range("O18").Value = 1.75
$ws.Range("AB21").Value = 15
$ws.Range("AD21").Value = 11.75
$ws.Range("H21").Value = 3.6
$ws.Range("I21").Value = 3.9
$ws.Range("R21").Value = 1.72
$ws.Range("Y21").Value = 26
$ws.Range("AF23").Value = 8
$ws.Range("G23").Value = 5.8
$ws.Range("H23").Value = 4
$ws.Range("W23").Value = 120
$ws.Range("X23").Value = 65
$ws.Range("Z23").Value = 11
$ws.Range("AH30").Value = 26
$ws.Range("AJ30").Value = 101
$ws.Range("G30").Value = 1.9
$ws.Range("I30").Value = 3.75
$ws.Range("R30").Value = 1.5
$ws.Range("S30").Value = 2.5
$ws.Range("T30").Value = 11
$ws.Range("V30").Value = 9
$ws.Range("Z30").Value = 17
$ws.Range("AB31").Value = 13
$ws.Range("G31").Value = 2.2
$ws.Range("H31").Value = 3.4
$ws.Range("I31").Value = 3.2
$ws.Range("J31").Value = 1.04
$ws.Range("K31").Value = 12
$ws.Range("N31").Value = 1.8
$ws.Range("O31").Value = 2
$ws.Range("V31").Value = 9
$ws.Range("K32").Value = 15
$ws.Range("N32").Value = 1.62
$ws.Range("O32").Value = 2.25
$ws.Range("AD33").Value = 12
$ws.Range("AE33").Value = 15
$ws.Range("AF33").Value = 10
$ws.Range("AI33").Value = 23
$ws.Range("G33").Value = 2.8
$ws.Range("I33").Value = 2.45
$ws.Range("Y33").Value = 23
$ws.Range("N35").Value = 1.36
$ws.Range("O35").Value = 3.1
